$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is purely numeric-looking (e.g. "1.00", "224.51")
# must have their number format forced to Text first, otherwise the COM
# layer auto-converts the assigned string into a real number (dropping the
# trailing zero / turning "1.00" into 1). Cells like "34.440.85" (two dots)
# are never parsed as numbers, so they do not need this treatment.

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Row-level value updates (price and 1h volume % changes) ---
$ws.Range("D2").Value = "34.440.85"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.799.52"
$ws.Range("E3").Value = "  +0.32%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.32%  "
Set-TextValue "D5" "224.51"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.36%  "
Set-TextValue "D8" "40.84"
$ws.Range("E8").Value = "  +13.12%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -1.05%  "
Set-TextValue "D11" "0.100"
$ws.Range("E11").Value = "  +4.16%  "
$ws.Range("D12").Value = "2.058.27"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.802.02"
$ws.Range("E13").Value = "  +0.01%  "
Set-TextValue "D14" "10.83"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "34.400.65"
$ws.Range("E15").Value = "  +0.37%  "
Set-TextValue "D16" "0.624"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  +0.31%  "
Set-TextValue "D18" "67.20"
$ws.Range("E18").Value = "  -1.98%  "
Set-TextValue "D19" "239.56"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "0.0₃0764"
$ws.Range("E20").Value = "  -0.60%  "
Set-TextValue "D21" "11.09"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("E22").Value = "  +0.39%  "
Set-TextValue "D23" "4.11"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -0.85%  "
Set-TextValue "D25" "172.06"
$ws.Range("E25").Value = "  +0.95%  "
Set-TextValue "D26" "7.64"
$ws.Range("E26").Value = "  -4.42%  "
Set-TextValue "D27" "17.31"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +0.56%  "
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  -0.01%  "
Set-TextValue "D31" "3.76"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -0.72%  "
Set-TextValue "D33" "0.0510"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "1.319.15"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("E44").Value = "  +0.24%  "
Set-TextValue "D45" "0.937"
$ws.Range("E45").Value = "  +0.71%  "
Set-TextValue "D46" "0.0519"
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("D47").Value = "1.959.43"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +1.45%  "
Set-TextValue "D49" "1.01"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  +1.51%  "

# --- Rows 37/38 swap: TrustWalletToken <-> Aave with new data ---
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D37" "85.92"
$ws.Range("E37").Value = "  +6.97%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.05"
$ws.Range("E38").Value = "  +0.44%  "

# --- Rows 40/41 swap: RenderToken <-> InjectiveProtocol with new data ---
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D40" "14.71"
$ws.Range("E40").Value = "  +12.22%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D41" "2.33"
$ws.Range("E41").Value = "  -0.16%  "
